# The "low mooring" run log previously had a leading placeholder row
# (fullpanel-nowind----depth580-mstop10-run2.csv) whose processed columns
# were never filled in. A rerun of load_or_update (now forcing a recompute)
# dropped that stale row, so the two fully-populated rows that followed it
# shift up to become rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
